# Update "想去人数" (interest count) values in column F for the "展览"
# and "全部类型" worksheets, as produced by the latest data scrape.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1894
$ws1.Range("F6").Value = 2686
$ws1.Range("F7").Value = 185
$ws1.Range("F8").Value = 95
$ws1.Range("F10").Value = 1565
$ws1.Range("F11").Value = 546
$ws1.Range("F16").Value = 181
$ws1.Range("F17").Value = 4
$ws1.Range("F19").Value = 225
$ws1.Range("F20").Value = 217
$ws1.Range("F22").Value = 210
$ws1.Range("F23").Value = 68
$ws1.Range("F24").Value = 1728
$ws1.Range("F25").Value = 40
$ws1.Range("F26").Value = 418
$ws1.Range("F27").Value = 63
$ws1.Range("F29").Value = 214
$ws1.Range("F31").Value = 441

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1894
$ws4.Range("F7").Value = 2686
$ws4.Range("F8").Value = 185
$ws4.Range("F9").Value = 95
$ws4.Range("F11").Value = 1565
$ws4.Range("F12").Value = 546
$ws4.Range("F17").Value = 181
$ws4.Range("F18").Value = 4
$ws4.Range("F20").Value = 225
$ws4.Range("F21").Value = 217
$ws4.Range("F23").Value = 210
$ws4.Range("F24").Value = 68
$ws4.Range("F25").Value = 1728
$ws4.Range("F26").Value = 40
$ws4.Range("F27").Value = 418
$ws4.Range("F28").Value = 63
$ws4.Range("F30").Value = 214
$ws4.Range("F32").Value = 441

$wb.Save()
